# Update "想去人数" (want-to-go count) figures in each sheet to match the
# latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13210
$ws1.Range("F6").Value  = 108
$ws1.Range("F7").Value  = 110
$ws1.Range("F8").Value  = 59
$ws1.Range("F10").Value = 24
$ws1.Range("F11").Value = 13153
$ws1.Range("F13").Value = 572
$ws1.Range("F14").Value = 8823
$ws1.Range("F15").Value = 7895
$ws1.Range("F16").Value = 226
$ws1.Range("F17").Value = 135
$ws1.Range("F20").Value = 1
$ws1.Range("F21").Value = 8
$ws1.Range("F27").Value = 82

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 24

# --- Sheet "全部类型" (all types combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 13210
$ws4.Range("F7").Value  = 108
$ws4.Range("F8").Value  = 110
$ws4.Range("F9").Value  = 59
$ws4.Range("F11").Value = 24
$ws4.Range("F12").Value = 13153
$ws4.Range("F14").Value = 572
$ws4.Range("F15").Value = 8823
$ws4.Range("F16").Value = 7895
$ws4.Range("F17").Value = 226
$ws4.Range("F18").Value = 135
$ws4.Range("F21").Value = 1
$ws4.Range("F22").Value = 8
$ws4.Range("F26").Value = 24
$ws4.Range("F30").Value = 82
